$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 69-154 down by two rows (to 71-156),
# matching the new weekly records inserted at the top of the series.
$src = $ws.Range("A69:R154")
$src.Copy($ws.Range("A71"))

# Row 69: new weekly record
$ws.Cells.Item(69, 1).Value = 10
$ws.Cells.Item(69, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(69, 3).Value = "La Araucanía"
$ws.Cells.Item(69, 4).Value = 44413
$ws.Cells.Item(69, 5).Value = 9
$ws.Cells.Item(69, 6).Value = 100112044
$ws.Cells.Item(69, 7).Value = "Perejil"
$ws.Cells.Item(69, 8).Value = "Sin especificar"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 40
$ws.Cells.Item(69, 11).Value = 3500
$ws.Cells.Item(69, 12).Value = 3500
$ws.Cells.Item(69, 13).Value = 3500
$ws.Cells.Item(69, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(69, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(69, 16).Value = 1167
$ws.Cells.Item(69, 17).Value = 3
$ws.Cells.Item(69, 18).Value = "Hortaliza"

# Row 70: new weekly record
$ws.Cells.Item(70, 1).Value = 10
$ws.Cells.Item(70, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(70, 3).Value = "La Araucanía"
$ws.Cells.Item(70, 4).Value = 44413
$ws.Cells.Item(70, 5).Value = 9
$ws.Cells.Item(70, 6).Value = 100112044
$ws.Cells.Item(70, 7).Value = "Perejil"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 50
$ws.Cells.Item(70, 11).Value = 3300
$ws.Cells.Item(70, 12).Value = 3300
$ws.Cells.Item(70, 13).Value = 3300
$ws.Cells.Item(70, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(70, 15).Value = "Región Metropolitana"
$ws.Cells.Item(70, 16).Value = 1100
$ws.Cells.Item(70, 17).Value = 3
$ws.Cells.Item(70, 18).Value = "Hortaliza"
